$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 18480.834
$ws.Range("J51").Value = 18480.834
$ws.Range("L51").Value = 18480.834
$ws.Range("N51").Value = -19448.834
$ws.Range("H55").Value = 189.375
$ws.Range("J55").Value = 324.625
$ws.Range("L55").Value = 324.625
$ws.Range("N55").Value = -752.625
$ws.Range("H116").Value = 7940226.5
$ws.Range("I116").Value = 10103797
$ws.Range("J116").Value = 7133
$ws.Range("K116").Value = 10103797
$ws.Range("L116").Value = 7133
$ws.Range("M116").Value = -10100355
$ws.Range("N116").Value = -14017
$ws.Range("H132").Value = 334123.06
$ws.Range("I132").Value = 369622.94
$ws.Range("K132").Value = 1108868.82
$ws.Range("M132").Value = -1106338.82
$ws.Range("H137").Value = 6216.7085
$ws.Range("I137").Value = 8676.454
$ws.Range("J137").Value = 4135.385
$ws.Range("K137").Value = 26029.362
$ws.Range("L137").Value = 12406.155
$ws.Range("M137").Value = -23479.362
$ws.Range("N137").Value = -17506.155
$ws.Range("H138").Value = 3990.0264
$ws.Range("I138").Value = 1753
$ws.Range("K138").Value = 5259
$ws.Range("M138").Value = -119

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 459.14285
$ws.Range("I5").Value = 459.14285
$ws.Range("K5").Value = 459.14285
$ws.Range("M5").Value = -347.14285
$ws.Range("H11").Value = 25001524
$ws.Range("I11").Value = 25001524
$ws.Range("K11").Value = 25001524
$ws.Range("M11").Value = -25001380
$ws.Range("H32").Value = 7939064
$ws.Range("I32").Value = 7353954.5
$ws.Range("K32").Value = 7353954.5
$ws.Range("M32").Value = -7353667.5
$ws.Range("H45").Value = 992
$ws.Range("I45").Value = 924.1667
$ws.Range("K45").Value = 924.1667
$ws.Range("M45").Value = -547.1667

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 459.14285
$ws.Range("I4").Value = 459.14285
$ws.Range("K4").Value = 459.14285
$ws.Range("M4").Value = -344.14285
$ws.Range("H20").Value = 6212218.5
$ws.Range("I20").Value = 14286276
$ws.Range("K20").Value = 14286276
$ws.Range("M20").Value = -14286029
$ws.Range("H105").Value = 4300
$ws.Range("I105").Value = 4300
$ws.Range("K105").Value = 4300
$ws.Range("M105").Value = -2553
$ws.Range("H107").Value = 2503576.8
$ws.Range("I107").Value = 3127847.8
$ws.Range("K107").Value = 3127847.8
$ws.Range("M107").Value = -3125927.8
$ws.Range("H134").Value = 3857207.2
$ws.Range("I134").Value = 5010119.5
$ws.Range("K134").Value = 15030358.5
$ws.Range("M134").Value = -15027823.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 8621.25
$ws.Range("I7").Value = 9810
$ws.Range("K7").Value = 9810
$ws.Range("M7").Value = -9697
$ws.Range("H16").Value = 19234754
$ws.Range("I16").Value = 35717304
$ws.Range("K16").Value = 35717304
$ws.Range("M16").Value = -35717017
$ws.Range("H22").Value = 3573866.5
$ws.Range("I22").Value = 4464583
$ws.Range("J22").Value = 11002
$ws.Range("K22").Value = 4464583
$ws.Range("L22").Value = 11002
$ws.Range("M22").Value = -4464233
$ws.Range("N22").Value = -11702
$ws.Range("H113").Value = 19234754
$ws.Range("I113").Value = 35717304
$ws.Range("K113").Value = 35717304
$ws.Range("M113").Value = -35715134
$ws.Range("H132").Value = 6438.229
$ws.Range("I132").Value = 4206.122
$ws.Range("J132").Value = 19512
$ws.Range("K132").Value = 12618.366
$ws.Range("L132").Value = 58536
$ws.Range("M132").Value = -10088.366
$ws.Range("N132").Value = -63596
$ws.Range("H141").Value = 713333.3
$ws.Range("J141").Value = 713333.3
$ws.Range("L141").Value = 713333.3
$ws.Range("N141").Value = -723693.3

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 24002.416
$ws.Range("I3").Value = 15007.25
$ws.Range("K3").Value = 45021.75
$ws.Range("M3").Value = -44909.75
$ws.Range("H114").Value = 1404.9231
$ws.Range("I114").Value = 1093.625
$ws.Range("J114").Value = 1903
$ws.Range("K114").Value = 3280.875
$ws.Range("L114").Value = 5709
$ws.Range("M114").Value = -26.875
$ws.Range("N114").Value = -12217
$ws.Range("H133").Value = 9640.888999999999
$ws.Range("I133").Value = 8025
$ws.Range("K133").Value = 24075
$ws.Range("M133").Value = -19015
$ws.Range("H138").Value = 2750
$ws.Range("I138").Value = 2750
$ws.Range("K138").Value = 8250
$ws.Range("M138").Value = -3110
$ws.Range("H139").Value = 3494.6924
$ws.Range("I139").Value = 2273.1667
$ws.Range("K139").Value = 6819.500100000001
$ws.Range("M139").Value = -1679.500100000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 11421306
$ws.Range("I11").Value = 12526631
$ws.Range("J11").Value = 7000004
$ws.Range("K11").Value = 12526631
$ws.Range("L11").Value = 7000004
$ws.Range("M11").Value = -12526492
$ws.Range("N11").Value = -7000282
$ws.Range("H12").Value = 10004
$ws.Range("J12").Value = 10004
$ws.Range("L12").Value = 10004
$ws.Range("N12").Value = -10284
$ws.Range("H14").Value = 11801447
$ws.Range("J14").Value = 950
$ws.Range("L14").Value = 950
$ws.Range("N14").Value = -1286
$ws.Range("H24").Value = 3335666.8
$ws.Range("I24").Value = 3335666.8
$ws.Range("K24").Value = 3335666.8
$ws.Range("M24").Value = -3335493.8
$ws.Range("H49").Value = 29334
$ws.Range("J49").Value = 29334
$ws.Range("L49").Value = 29334
$ws.Range("N49").Value = -29702
$ws.Range("H126").Value = 13520066
$ws.Range("I126").Value = 20836886
$ws.Range("K126").Value = 62510658
$ws.Range("M126").Value = -62508188
$ws.Range("H132").Value = 66673376
$ws.Range("I132").Value = 90915820
$ws.Range("K132").Value = 272747460
$ws.Range("M132").Value = -272744930

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H7").Value = 8884.9
$ws.Range("I7").Value = 9413.916999999999
$ws.Range("K7").Value = 9413.916999999999
$ws.Range("M7").Value = -9301.916999999999
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H16").Value = 2748.45
$ws.Range("I16").Value = 595
$ws.Range("K16").Value = 595
$ws.Range("M16").Value = -425
$ws.Range("H40").Value = 5843.143
$ws.Range("I40").Value = 5812.2354
$ws.Range("J40").Value = 5974.5
$ws.Range("K40").Value = 5812.2354
$ws.Range("L40").Value = 5974.5
$ws.Range("M40").Value = -5676.2354
$ws.Range("N40").Value = -6246.5
$ws.Range("H61").Value = 5165.1284
$ws.Range("I61").Value = 3967
$ws.Range("K61").Value = 3967
$ws.Range("M61").Value = -3765
$ws.Range("H99").Value = 47499.75
$ws.Range("I99").Value = 47499.5
$ws.Range("K99").Value = 47499.5
$ws.Range("M99").Value = -44504.5
$ws.Range("H113").Value = 5165.1284
$ws.Range("I113").Value = 3967
$ws.Range("K113").Value = 3967
$ws.Range("M113").Value = -1797
$ws.Range("H122").Value = 6561.2
$ws.Range("I122").Value = 5666.3335
$ws.Range("K122").Value = 16999.0005
$ws.Range("M122").Value = -14549.0005
$ws.Range("H126").Value = 8884.9
$ws.Range("I126").Value = 9413.916999999999
$ws.Range("K126").Value = 28241.751
$ws.Range("M126").Value = -25771.751

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 6258213.5
$ws.Range("I9").Value = 8339116.5
$ws.Range("J9").Value = 15503.5
$ws.Range("K9").Value = 8339116.5
$ws.Range("L9").Value = 15503.5
$ws.Range("M9").Value = -8338976.5
$ws.Range("N9").Value = -15783.5
$ws.Range("H10").Value = 12500875
$ws.Range("I10").Value = 16667333
$ws.Range("J10").Value = 1500
$ws.Range("K10").Value = 16667333
$ws.Range("L10").Value = 1500
$ws.Range("M10").Value = -16667164
$ws.Range("N10").Value = -1838
$ws.Range("H14").Value = 528.7143
$ws.Range("I14").Value = 500.16666
$ws.Range("J14").Value = 700
$ws.Range("K14").Value = 500.16666
$ws.Range("L14").Value = 700
$ws.Range("M14").Value = -332.16666
$ws.Range("N14").Value = -1036
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H56").Value = 68500
$ws.Range("J56").Value = 67000
$ws.Range("L56").Value = 67000
$ws.Range("N56").Value = -68428
$ws.Range("H81").Value = 1766.1578
$ws.Range("I81").Value = 1647.3636
$ws.Range("J81").Value = 1929.5
$ws.Range("K81").Value = 3294.7272
$ws.Range("L81").Value = 3859
$ws.Range("M81").Value = -2233.7272
$ws.Range("N81").Value = -5981
$ws.Range("H84").Value = 1766.1578
$ws.Range("I84").Value = 1647.3636
$ws.Range("J84").Value = 1929.5
$ws.Range("K84").Value = 16473.636
$ws.Range("L84").Value = 19295
$ws.Range("M84").Value = -11169.636
$ws.Range("N84").Value = -29903
$ws.Range("H126").Value = 3804.8125
$ws.Range("I126").Value = 2516.6667
$ws.Range("K126").Value = 7550.000100000001
$ws.Range("M126").Value = -5080.000100000001
$ws.Range("H132").Value = 10892.412
$ws.Range("J132").Value = 12937.5
$ws.Range("L132").Value = 38812.5
$ws.Range("N132").Value = -43872.5

Write-Host "Applied all market-data updates"